$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.090.30'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '3.872.52'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Formula = "'599.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Formula = "'167.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("D7").Value = '3.869.61'
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").Formula = "'6.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Formula = "'0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Formula = "'36.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").Value = '4.518.08'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("D16").Value = '3.864.99'
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").Value = '68.111.20'
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").Formula = "'18.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.66%  '
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").Formula = "'10.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").Formula = "'466.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.28%  '
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("D26").Formula = "'2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.75%  '
$ws.Range("D27").Formula = "'12.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.31%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Formula = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Formula = "'10.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").Value = '4.022.23'
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").Formula = "'7.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.58%  '
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("E34").Value = '  -2.31%  '
$ws.Range("D35").Formula = "'9.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("D36").Value = '3.846.56'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D38").Formula = "'3.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.28%  '
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("D41").Formula = "'5.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("D44").Formula = "'434.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Formula = "'47.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").Formula = "'8.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Formula = "'0.000280"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.87%  '
$ws.Range("D50").Formula = "'40.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.42%  '
$ws.Range("D51").Formula = "'143.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.44%  '
